# "Added time into last_run, of workflows list"
#
# The last_run column (E) was storing only a date (serial number, formatted
# with a custom date/time numFmt). We now stamp each row's last_run with a
# full date+time text value, so the old numeric date-format styling is no
# longer needed on these cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the old date/time number formatting from the last_run cells - they
# will now hold plain date+time text instead of a formatted date serial.
$ws.Range("E2:E10").ClearFormats()

# Stamp each workflow row with its last_run date and time.
$ws.Range("E2").Value  = "2026-01-28 00:00:00"
$ws.Range("E3").Value  = "2026-01-28 23:18:10"
$ws.Range("E4").Value  = "2026-01-28 00:00:00"
$ws.Range("E5").Value  = "2026-01-28 00:00:00"
$ws.Range("E6").Value  = "2026-01-28 00:00:00"
$ws.Range("E7").Value  = "2026-01-28 00:00:00"
$ws.Range("E8").Value  = "2026-01-28 00:00:00"
$ws.Range("E9").Value  = "2026-01-28 00:00:00"
$ws.Range("E10").Value = "2026-01-28 00:00:00"
